# Quarterly update: drop the oldest period (6-month 1399/06 column) and
# append the newly published 12-month 1401/12 period; also correct a
# handful of values produced by the revised read_price algorithm.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the obsolete "6 ماهه منتهی به 1399/06" column. Excel shifts
#    everything from E..M left into D..L (and compacts the now-unused
#    shared strings for that period/date automatically).
$ws.Columns("D").Delete()

# 2) Clone column L's formatting into the new trailing column M so the
#    freshly appended period lines up visually with the rest of the table
#    (bounded to the used rows so Excel doesn't materialise 1M+ rows).
$ws.Range("L1:L28").Copy()
$ws.Range("M1:M28").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Columns("M").ColumnWidth = $ws.Columns("E").ColumnWidth

# 3) New period header + publish date for column M.
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"
$ws.Range("M9").Value = "1402-01-29"

# 4) New financial figures for the freshly published 12-month 1401/12 period.
$ws.Range("M11").Value = 410501446
$ws.Range("M12").Value = -219085957
$ws.Range("M13").Value = 191415489
$ws.Range("M14").Value = -9950238
$ws.Range("M15").Value = 0
$ws.Range("M16").Value = 0
$ws.Range("M17").Value = 181465251
$ws.Range("M18").Value = 0
$ws.Range("M19").Value = 23002071
$ws.Range("M20").Value = 204467322
$ws.Range("M21").Value = 0
$ws.Range("M22").Value = 204467322
$ws.Range("M23").Value = 0
$ws.Range("M24").Value = 204467322
$ws.Range("M25").Value = 1062
$ws.Range("M26").Value = 192500000
$ws.Range("M27").Value = 753

# 5) Corrections from the revised read_price algorithm applied across the
#    already-published periods (rows 19, 20, 22, 24 column I, and the
#    recomputed "EPS on latest capital" row 27 across all periods).
$ws.Range("I19").Value = 11181062
$ws.Range("I20").Value = 293203737
$ws.Range("I22").Value = 269476894
$ws.Range("I24").Value = 269476894

$ws.Range("D27").Value = 348
$ws.Range("E27").Value = 630
$ws.Range("F27").Value = 221
$ws.Range("G27").Value = 475
$ws.Range("H27").Value = 672
$ws.Range("I27").Value = 993
$ws.Range("J27").Value = 254
$ws.Range("K27").Value = 354
$ws.Range("L27").Value = 516
